$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row for "Canada" right before "Cargolux" (original row 9).
# ---------------------------------------------------------------------------
$ws.Rows("9:9").Insert()
$ws.Range("A9").Value = "Canada"
$ws.Range("B9").Value = "014"
$ws.Range("C9").Value = "in progress"
$ws.Range("D9").Value = "https://www.aircanada.com/cargo/en/tools-forms/"

# British (row 8) is no longer "in progress".
$ws.Range("C8").ClearContents()

# ---------------------------------------------------------------------------
# 2) Insert two new rows for "Swiss" and "Tampa" right before "Turkish".
#    After step 1, Turkish (originally row 24) now sits at row 25.
# ---------------------------------------------------------------------------
$ws.Rows("26:27").Insert()

$ws.Range("A26").Value = "Swiss"
$ws.Range("B26").Value = "724"
$ws.Range("C26").Value = "in progress"
$ws.Range("D26").Value = "https://www.swissworldcargo.com/track_n_trace"

$ws.Range("A27").Value = "Tampa"
$ws.Range("B27").Value = "729"
$ws.Range("C27").Value = "in progress"
$ws.Range("D27").Value = "http://www.aviancacargo.com/index.aspx"

# ---------------------------------------------------------------------------
# 3) Rebuild the hyperlinks so every URL cell in column D is a live link
#    again (row-insert does not automatically move hyperlink anchors).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D2"), "https://www.airbridgecargo.com/en/tracking/")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://www.airtahitinui.com/us-en/online-cargo-tracking")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://www.afklcargo.com/WW/en/local/app/index.jsp", "/tntsinglesearch")
$ws.Hyperlinks.Add($ws.Range("D11"), "http://www.cathaypacificcargo.com/ManageYourShipment/TrackYourShipment/tabid/108/SingleAWBNo/160-05480334-/language/en-US/Default.aspx")
$ws.Hyperlinks.Add($ws.Range("D12"), "https://cargo.china-airlines.com/CCNetv2/content/manage/ShipmentTracking.aspx?")
$ws.Hyperlinks.Add($ws.Range("D13"), "https://aviationcargo.dhl.com/aviationcargo/track/")
$ws.Hyperlinks.Add($ws.Range("D14"), "https://skychain.emirates.com/skychain/app?service=page/nwp:Trackshipmt&initial=y")
$ws.Hyperlinks.Add($ws.Range("D15"), "http://www.brcargo.com/ec_web/Default.aspx?Parm2=191&Parm3=undefined")
$ws.Hyperlinks.Add($ws.Range("D16"), "http://www.jal.co.jp/en/jalcargo/inter/awb/")
$ws.Hyperlinks.Add($ws.Range("D18"), "https://www.afklcargo.com/WW/en/local/app/index.jsp", "/tntsinglesearch")
$ws.Hyperlinks.Add($ws.Range("D17"), "https://lufthansa-cargo.com/eservices/etracking")
$ws.Hyperlinks.Add($ws.Range("D19"), "https://cargo.koreanair.com/en/tracking?")
$ws.Hyperlinks.Add($ws.Range("D20"), "http://www.maskargo.com/online_awb_info/index.php")
$ws.Hyperlinks.Add($ws.Range("D22"), "https://freight.qantas.com/online-tracking.html?")
$ws.Hyperlinks.Add($ws.Range("D23"), "http://www.qrcargo.com/trackshipment")
$ws.Hyperlinks.Add($ws.Range("D24"), "http://www.siacargo.com/ccn/ShipmentTrack.aspx")
$ws.Hyperlinks.Add($ws.Range("D25"), "https://www.skyteam.com/en/cargo/track-shipment/")
$ws.Hyperlinks.Add($ws.Range("D29"), "https://www.unitedcargo.com/OurNetwork/TrackingCargo1512/Tracking.jsp")
$ws.Hyperlinks.Add($ws.Range("D28"), "https://www.turkishcargo.com.tr/en/online-services/shipment-tracking")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://www.aacargo.com/AACargo/tracking")
$ws.Hyperlinks.Add($ws.Range("D6"), "https://mycargo.amerijet.com/tracking")
$ws.Hyperlinks.Add($ws.Range("D21"), "https://www.anacargo.jp/en/int/")
$ws.Hyperlinks.Add($ws.Range("D7"), "https://www.asianacargo.com/tracking/viewTraceAirWaybill.do?lang=en")
$ws.Hyperlinks.Add($ws.Range("D8"), "https://www.iagcargo.com/en/home")
$ws.Hyperlinks.Add($ws.Range("D30"), "https://cargo.virgin-atlantic.com/gb/en/track/track-your-cargo.html?prefix=932&number=56409673&track=go")
$ws.Hyperlinks.Add($ws.Range("D10"), "https://cvtnt.champ.aero/trackntrace")
$ws.Hyperlinks.Add($ws.Range("D9"), "https://www.aircanada.com/cargo/en/tools-forms/")
$ws.Hyperlinks.Add($ws.Range("D26"), "https://www.swissworldcargo.com/track_n_trace")
$ws.Hyperlinks.Add($ws.Range("D27"), "http://www.aviancacargo.com/index.aspx")

# ---------------------------------------------------------------------------
# 4) Selection / active cell bookkeeping to match the final workbook state.
# ---------------------------------------------------------------------------
$ws.Range("B8:B30").Select()
